{"js": "// Replace the 25 \"three-digit \u00f7 one-digit\" answer strings in the table\n// cells with their new values, keeping every other part of the document\n// (formatting, empty rows, the date paragraph, etc.) untouched.\nconst replacements = [\n  [\"808\u00f77=115, 3\", \"303\u00f73=101, 0\"],\n  [\"515\u00f72=257, 1\", \"433\u00f75=86, 3\"],\n  [\"173\u00f77=24, 5\", \"342\u00f74=85, 2\"],\n  [\"418\u00f73=139, 1\", \"357\u00f79=39, 6\"],\n  [\"798\u00f79=88, 6\", \"848\u00f78=106, 0\"],\n  [\"547\u00f77=78, 1\", \"590\u00f79=65, 5\"],\n  [\"187\u00f77=26, 5\", \"898\u00f76=149, 4\"],\n  [\"351\u00f73=117, 0\", \"206\u00f78=25, 6\"],\n  [\"542\u00f75=108, 2\", \"332\u00f78=41, 4\"],\n  [\"221\u00f76=36, 5\", \"743\u00f75=148, 3\"],\n  [\"575\u00f74=143, 3\", \"485\u00f74=121, 1\"],\n  [\"356\u00f75=71, 1\", \"367\u00f75=73, 2\"],\n  [\"807\u00f72=403, 1\", \"748\u00f72=374, 0\"],\n  [\"752\u00f73=250, 2\", \"476\u00f72=238, 0\"],\n  [\"673\u00f72=336, 1\", \"867\u00f78=108, 3\"],\n  [\"761\u00f75=152, 1\", \"834\u00f79=92, 6\"],\n  [\"239\u00f76=39, 5\", \"623\u00f73=207, 2\"],\n  [\"790\u00f75=158, 0\", \"483\u00f78=60, 3\"],\n  [\"409\u00f78=51, 1\", \"740\u00f79=82, 2\"],\n  [\"988\u00f79=109, 7\", \"450\u00f77=64, 2\"],\n  [\"454\u00f73=151, 1\", \"241\u00f78=30, 1\"],\n  [\"277\u00f74=69, 1\", \"108\u00f75=21, 3\"],\n  [\"167\u00f79=18, 5\", \"163\u00f77=23, 2\"],\n  [\"946\u00f78=118, 2\", \"540\u00f73=180, 0\"],\n  [\"321\u00f73=107, 0\", \"816\u00f78=102, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWildcards: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"three-digit \u00f7 one-digit\" answer strings in the table\n# cells with their new values, keeping every other part of the document\n# (formatting, empty rows, the date paragraph, etc.) untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"808\u00f77=115, 3\", \"303\u00f73=101, 0\"),\n  @(\"515\u00f72=257, 1\", \"433\u00f75=86, 3\"),\n  @(\"173\u00f77=24, 5\", \"342\u00f74=85, 2\"),\n  @(\"418\u00f73=139, 1\", \"357\u00f79=39, 6\"),\n  @(\"798\u00f79=88, 6\", \"848\u00f78=106, 0\"),\n  @(\"547\u00f77=78, 1\", \"590\u00f79=65, 5\"),\n  @(\"187\u00f77=26, 5\", \"898\u00f76=149, 4\"),\n  @(\"351\u00f73=117, 0\", \"206\u00f78=25, 6\"),\n  @(\"542\u00f75=108, 2\", \"332\u00f78=41, 4\"),\n  @(\"221\u00f76=36, 5\", \"743\u00f75=148, 3\"),\n  @(\"575\u00f74=143, 3\", \"485\u00f74=121, 1\"),\n  @(\"356\u00f75=71, 1\", \"367\u00f75=73, 2\"),\n  @(\"807\u00f72=403, 1\", \"748\u00f72=374, 0\"),\n  @(\"752\u00f73=250, 2\", \"476\u00f72=238, 0\"),\n  @(\"673\u00f72=336, 1\", \"867\u00f78=108, 3\"),\n  @(\"761\u00f75=152, 1\", \"834\u00f79=92, 6\"),\n  @(\"239\u00f76=39, 5\", \"623\u00f73=207, 2\"),\n  @(\"790\u00f75=158, 0\", \"483\u00f78=60, 3\"),\n  @(\"409\u00f78=51, 1\", \"740\u00f79=82, 2\"),\n  @(\"988\u00f79=109, 7\", \"450\u00f77=64, 2\"),\n  @(\"454\u00f73=151, 1\", \"241\u00f78=30, 1\"),\n  @(\"277\u00f74=69, 1\", \"108\u00f75=21, 3\"),\n  @(\"167\u00f79=18, 5\", \"163\u00f77=23, 2\"),\n  @(\"946\u00f78=118, 2\", \"540\u00f73=180, 0\"),\n  @(\"321\u00f73=107, 0\", \"816\u00f78=102, 0\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute(\n    $oldText,   # FindText\n    $true,      # MatchCase\n    $false,     # MatchWholeWord\n    $false,     # MatchWildcards\n    $false,     # MatchSoundsLike\n    $false,     # MatchAllWordForms\n    $true,      # Forward\n    1,          # Wrap (wdFindContinue)\n    $false,     # Format\n    $newText,   # ReplaceWith\n    2           # Replace (wdReplaceAll)\n  ) | Out-Null\n}\n"}
